# Fix shared-string values that had an accidental trailing space
# ("craft keys" i.e. the English category labels lost their trailing
# whitespace mistake) and move the active selection to A9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value2 = "Work"
$ws.Range("A5").Value2 = "Shopping"
$ws.Range("A8").Value2 = "Accompaniment"

$ws.Range("A9").Select()
